# "filter and sort feature"
#
# Renames the "Customer" sheet to "Filter", rewrites the Login-sheet
# credential row (and re-links its hyperlink/style), and populates the
# Filter sheet with the Employee/Card filter table. Also makes the
# Filter sheet the active tab, matching the saved workbook view.

$wb = $excel.ActiveWorkbook

# --- Sheet1 "Login": update the credential values ---------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A2").Value = "FPTTeam"
$ws1.Range("B2").Value = "123@a123"

# B2 becomes a hyperlink (mirroring A2's existing mailto link) and picks
# up the built-in "Hyperlink" cell style, same as A2.
$ws1.Hyperlinks.Add($ws1.Range("B2"), "mailto:123@a123")
$ws1.Range("B2").Style = "Hyperlink"

# --- Sheet2 "Customer" -> "Filter": rename + add data ------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Filter"

$ws2.Range("A1").Value = "EmployeeName"
$ws2.Range("B1").Value = "EmployeeId"
$ws2.Range("C1").Value = "Departments"
$ws2.Range("D1").Value = "CardStatus"
$ws2.Range("E1").Value = "CardType"

$ws2.Range("A2").Value = "mmaa"
$ws2.Range("B2").Value = "A4352"
$ws2.Range("C2").Value = "Admin"
$ws2.Range("D2").Value = "Activated"
$ws2.Range("E2").Value = "Instant Issue"

# Column widths roughly matching the bestFit results captured in the
# saved workbook.
$ws2.Columns.Item(1).ColumnWidth = 12.996651785714286
$ws2.Columns.Item(2).ColumnWidth = 9.711495535714286
$ws2.Columns.Item(3).ColumnWidth = 10.285714285714286
$ws2.Columns.Item(4).ColumnWidth = 9.141183035714286
$ws2.Columns.Item(5).ColumnWidth = 8.141183035714286

# Filter is now the visible/active sheet, zoomed to 130%, with the same
# selection (H12) that the Login sheet had.
[void]$ws2.Range("H12").Select()
$excel.ActiveWindow.Zoom = 130
